$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new row 20: access page for orders management
$ws.Range("A20").Value = "accès à la page de gestion des commande"
$ws.Range("C20").Value = "x"
$ws.Range("D20").Value = "x"

# Match formatting (center/center alignment) used by rows 17-19
$ws.Range("C20:D20").HorizontalAlignment = -4108
$ws.Range("C20:D20").VerticalAlignment = -4108

# Update selection to match diff
$ws.Range("D23").Select()
